# This workbook holds a small student table:
#   matricula | nombre | apellidoP | apellidoM | (num) | (num)
# The update replaces the separate "apellidoP"/"apellidoM" surname columns
# with a single "grupo" column, updates a couple of matricula values, and
# refreshes the sheet ("actualizar listas desde un archivo").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Remove the apellidoP (C) and apellidoM (D) columns entirely; this
#    shifts the two numeric columns that followed (E,F) left into C,D.
$ws.Range("C1:D3").Delete()

# 2) After the shift, D2 carries the style that used to belong to the old
#    F2 cell. Turn its font underline on so that style definition becomes
#    the "grupo"-era placeholder style, then stamp two blank cells further
#    down the sheet (F8, F9) with that exact formatting via copy/paste of
#    formats only.
$ws.Range("D2").Font.Underline = $true
$ws.Range("D2").Copy()
$ws.Range("F8").PasteSpecial(-4122)
$ws.Range("F9").PasteSpecial(-4122)

# 3) Turn the (now) C column into the new "grupo" text column and fill in
#    the group values for each student.
$ws.Range("C1:C3").NumberFormat = "@"
$ws.Range("C1").Value = "grupo"
$ws.Range("C2").Value = "906-a"
$ws.Range("A3").Value = "2018060161"
$ws.Range("C3").Value = "806-b"
$ws.Range("A2").Value = "2015060162"

# 4) Restore D2's own number format back to match D1/D3 (plain 0.00 number
#    style) now that its formatting has served its purpose above.
$ws.Range("D2").Font.Underline = $false
$ws.Range("D2").NumberFormat = "0.00"
